$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C values (rows 4-24)
$ws.Range("C4").Value = "Üzleti folyamat feltérképezés, Felhasználókezelés, kosár kezelése, utánanézni, mik a lehetőségek"
$ws.Range("C5").Value = "Felhasználókezelés, kosár kezelése - backend megvalósítás, backend - NRT hozzáadása, Domain fejlesztése, refaktorálás"
$ws.Range("C6").Value = "Felhasználókezelés backenden - NRT miatt változtatások"
$ws.Range("C7").Value = "Felhasználókezelés backenden - refresh token megoldása"
$ws.Range("C8").Value = "Felhasználókezelés backenden - refresh token megoldása"
$ws.Range("C9").Value = "Kliens alkalmazás létrehozása, függőségek telepítése, Angular Material téma kialakítása"
$ws.Range("C10").Value = "Kliens alkalmazás basic layout, Felhasználókezelés - refresh token autentikáció"
$ws.Range("C11").Value = "RxJS betanulás - Subject, BehaviorSubject, Operators (pipe, map, tap, switchMap...), Kliens alkalmazás reszponzív header"
$ws.Range("C12").Value = "Kosár kezelése - CartService megvalósítása"
$ws.Range("C13").Value = "Felhasználókezelés, kosár kezelése - felület létrehozás - bejelentkezés, profil, rendeléseim képernyő (tipikus user flow)"
$ws.Range("C14").Value = "Shop filters - utánaolvasni, hogyan szokás elkészíteni, backend queryk? Megvalósítása, webes kliensben megvalósítás - shopban filter felület, filter service?"
$ws.Range("C15").Value = "Admin felület (áruk hozzáadása, törlése, szerkesztése, rendelések nézése)"
$ws.Range("C16").Value = "Admin felület (áruk hozzáadása, törlése, szerkesztése, rendelések nézése)"
$ws.Range("C17").Value = "Fizetés - utánanézni, milyen lehetőségek vannak, hogyan lehet beépíteni őket, ki lehet-e próbálni őket"
$ws.Range("C18").Value = "Webes fizetés megvalósítása"
$ws.Range("C19").Value = "Kereső optimalizálás - utánanézni, milyen módszerek vannak rá, hogyan érdemes csinálni, implementálni"
$ws.Range("C20").Value = "Android kliens"
$ws.Range("C21").Value = "Android kliens"
$ws.Range("C22").Value = "Android kliens"
$ws.Range("C23").Value = "Android kliens"
$ws.Range("C24").Value = "Ajánló motor, tesztelés, deployment"

# Set column D values
$ws.Range("D6").Value = "csúszás"
$ws.Range("D7").Value = "csúszás"
$ws.Range("D8").Value = "csúszás"
$ws.Range("D15").Value = "Dokumentáció írás (technológiák bemutatása, végleges részek bemutatása)"

# Set column D width (closest achievable value to target 77.140625 character-width units,
# which the engine internally quantizes to 1/6-character pixel steps)
$ws.Columns("D").ColumnWidth = 76.27

# Update selection
$ws.Range("D15").Select()